# Refresh the cryptocurrency price / volume(1h) table with the latest
# scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Addr, $Text)
    $range = $ws.Range($Addr)
    if ($Text -match "^-?\d+(\.\d+)?$") {
        # Excel would otherwise auto-parse this literal as a number
        # (e.g. "19.55" or "1.00"). Force literal text via the classic
        # leading-apostrophe entry convention, then strip the implicit
        # "Text" number-format style that leaves behind so the cell
        # keeps its original (default) formatting.
        $range.Value = "'" + $Text
        $range.Style = "Normal"
    } else {
        $range.Value = $Text
    }
}

# Row 2
Set-TextValue 'D2' '26.620.34'
Set-TextValue 'E2' '  -1.80%  '

# Row 3
Set-TextValue 'D3' '1.587.73'
Set-TextValue 'E3' '  -2.26%  '

# Row 4
Set-TextValue 'E4' '  +0.08%  '

# Row 5
Set-TextValue 'D5' '210.83'
Set-TextValue 'E5' '  -1.94%  '

# Row 6
Set-TextValue 'E6' '  -2.44%  '

# Row 7
Set-TextValue 'E7' '  +0.06%  '

# Row 8
Set-TextValue 'E8' '  -2.29%  '

# Row 9
Set-TextValue 'E9' '  -1.91%  '

# Row 10
Set-TextValue 'D10' '19.55'

# Row 11
Set-TextValue 'E11' '  -1.60%  '

# Row 12
Set-TextValue 'D12' '1.809.91'
Set-TextValue 'E12' '  -2.27%  '

# Row 13
Set-TextValue 'D13' '1.595.80'
Set-TextValue 'E13' '  -1.39%  '

# Row 14
Set-TextValue 'E14' '  -2.80%  '

# Row 15
Set-TextValue 'E15' '  -4.06%  '

# Row 16
Set-TextValue 'D16' '64.69'
Set-TextValue 'E16' '  +0.01%  '

# Row 17
Set-TextValue 'D17' '26.595.68'

# Row 18
Set-TextValue 'E18' '  -2.10%  '

# Row 19
Set-TextValue 'B19' 'BitcoinCash'
Set-TextValue 'C19' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D19' '208.38'
Set-TextValue 'E19' '  -3.89%  '

# Row 20
Set-TextValue 'B20' 'Dai'
Set-TextValue 'C20' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D20' '1.00'
Set-TextValue 'E20' '  +0.01%  '

# Row 21
Set-TextValue 'D21' '6.71'
Set-TextValue 'E21' '  -3.45%  '

# Row 22
Set-TextValue 'E22' '  -2.74%  '

# Row 23
Set-TextValue 'D23' '2.34'
Set-TextValue 'E23' '  -3.11%  '

# Row 24
Set-TextValue 'D24' '8.86'
Set-TextValue 'E24' '  -2.23%  '

# Row 25
Set-TextValue 'D25' '146.88'
Set-TextValue 'E25' '  -0.73%  '

# Row 26
Set-TextValue 'E26' '  +0.11%  '

# Row 27
Set-TextValue 'E27' '  -0.73%  '

# Row 28
Set-TextValue 'E28' '  -3.21%  '

# Row 29
Set-TextValue 'D29' '15.26'
Set-TextValue 'E29' '  -2.17%  '

# Row 30
Set-TextValue 'E30' '  -0.06%  '

# Row 31
Set-TextValue 'E31' '  -2.01%  '

# Row 32
Set-TextValue 'D32' '3.23'
Set-TextValue 'E32' '  -3.90%  '

# Row 33
Set-TextValue 'D33' '0.679'
Set-TextValue 'E33' '  +22.81%  '

# Row 34
Set-TextValue 'D34' '2.90'
Set-TextValue 'E34' '  -3.04%  '

# Row 35
Set-TextValue 'D35' '1.306.95'
Set-TextValue 'E35' '  -3.00%  '

# Row 36
Set-TextValue 'D36' '2.43'
Set-TextValue 'E36' '  -0.98%  '

# Row 37
Set-TextValue 'E37' '  -5.64%  '

# Row 38
Set-TextValue 'E38' '  -3.20%  '

# Row 39
Set-TextValue 'E39' '  -3.38%  '

# Row 40
Set-TextValue 'E40' '  +0.09%  '

# Row 41
Set-TextValue 'D41' '0.791'
Set-TextValue 'E41' '  -1.45%  '

# Row 42
Set-TextValue 'E42' '  +2.59%  '

# Row 43
Set-TextValue 'E43' '  -2.79%  '

# Row 44
Set-TextValue 'D44' '62.65'

# Row 45
Set-TextValue 'D45' '1.723.26'
Set-TextValue 'E45' '  -2.11%  '

# Row 46
Set-TextValue 'D46' '89.47'
Set-TextValue 'E46' '  -1.34%  '

# Row 47
Set-TextValue 'E47' '  -1.04%  '

# Row 48
Set-TextValue 'D48' '0.837'
Set-TextValue 'E48' '  -3.07%  '

# Row 49
Set-TextValue 'D49' '0.0504'
Set-TextValue 'E49' '  -1.78%  '

# Row 50
Set-TextValue 'D50' '0.0978'
Set-TextValue 'E50' '  -1.60%  '

# Row 51
Set-TextValue 'D51' '7.49'
Set-TextValue 'E51' '  -1.79%  '

